$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.460.25'
$ws.Range("D3").Value = '3.688.74'
$ws.Range("E3").Value = '  -3.20%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'687.60"
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").Value = "'161.71"
$ws.Range("E6").Value = '  -5.47%  '
$ws.Range("D7").Value = '3.687.53'
$ws.Range("E7").Value = '  -3.14%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -5.53%  '
$ws.Range("E10").Value = '  -8.31%  '
$ws.Range("D11").Value = "'7.36"
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").Value = "'0.438"
$ws.Range("E12").Value = '  -9.84%  '
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = '  -6.05%  '
$ws.Range("D14").Value = '4.312.32'
$ws.Range("E14").Value = '  -3.19%  '
$ws.Range("D15").Value = "'32.98"
$ws.Range("E15").Value = '  -8.39%  '
$ws.Range("D16").Value = '3.688.34'
$ws.Range("E16").Value = '  -3.27%  '
$ws.Range("D17").Value = '69.435.78'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = "'16.02"
$ws.Range("E19").Value = '  -8.56%  '
$ws.Range("D20").Value = "'6.50"
$ws.Range("E20").Value = '  -10.11%  '
$ws.Range("D21").Value = "'475.60"
$ws.Range("E21").Value = '  -7.60%  '
$ws.Range("D22").Value = "'9.93"
$ws.Range("E22").Value = '  -5.10%  '
$ws.Range("D23").Value = "'0.656"
$ws.Range("E23").Value = '  -8.17%  '
$ws.Range("D24").Value = "'79.86"
$ws.Range("E24").Value = '  -5.02%  '
$ws.Range("D25").Value = '3.836.14'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("E26").Value = '  -9.28%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -10.12%  '
$ws.Range("D29").Value = "'9.25"
$ws.Range("E29").Value = '  -10.93%  '
$ws.Range("D30").Value = "'1.79"
$ws.Range("E30").Value = '  -11.59%  '
$ws.Range("E31").Value = '  -10.28%  '
$ws.Range("D32").Value = "'6.76"
$ws.Range("E32").Value = '  -7.96%  '
$ws.Range("D33").Value = "'2.04"
$ws.Range("E33").Value = '  -8.53%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'26.81"
$ws.Range("E35").Value = '  -7.95%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.165"
$ws.Range("E36").Value = '  -4.58%  '
$ws.Range("D37").Value = '3.655.79'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("D38").Value = "'8.32"
$ws.Range("E38").Value = '  -9.41%  '
$ws.Range("D39").Value = "'6.26"
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = '  -4.54%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = "'0.0916"
$ws.Range("E42").Value = '  -9.25%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = "'0.951"
$ws.Range("E44").Value = '  -6.35%  '
$ws.Range("D45").Value = "'163.34"
$ws.Range("E45").Value = '  -5.19%  '
$ws.Range("D46").Value = "'48.31"
$ws.Range("E46").Value = '  -3.07%  '
$ws.Range("D47").Value = "'30.24"
$ws.Range("E47").Value = '  +2.59%  '
$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = '  -15.60%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = "'1.32"
$ws.Range("E49").Value = '  -3.94%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = "'0.000281"
$ws.Range("E50").Value = '  -8.96%  '
$ws.Range("E51").Value = '  -3.08%  '
